# timesheet.xlsx update — "progress on gun animations"
# D column = "LO3 (VFX+SFX+animation)" hours; log an extra 22+10 minutes on
# the 45335 ("row 34") date, then leave the selection on D35 (next blank
# LO3 entry row), matching where the author clicked next.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34 (date 45335) previously had D34 = 0 (via shared formula). Give it
# its own explicit formula/value, which ripples through all the dependent
# totals/averages elsewhere on the sheet automatically.
$ws.Range("D34").Formula = "=(1/60)*(22+10)"

# Move the active selection from B35 to D35.
$ws.Range("D35").Select()
